# The edit appends " change" and " 2" (as two separate runs, mirroring
# two distinct typing/insert actions) to the end of the "asim2" paragraph,
# and moves the "_GoBack" bookmark (which Word auto-maintains at the
# location of the most recent edit) from the start of the paragraph to
# its new end.

$d = $word.ActiveDocument

# --- 1. Append the new text as two separate runs -------------------------
# Locate the end of the (only) paragraph's text, i.e. just before its
# paragraph mark, and insert the first new chunk there.
$para = $d.Paragraphs(1).Range
$insertPos = $para.End - 1
$r = $d.Range($insertPos, $insertPos)
$r.InsertAfter(" change")
$r.Collapse(0)

# Insert the final chunk, plus a one-character placeholder ("X") so that
# the eventual bookmark position is not the very last character in the
# paragraph (inserting a bookmark collapsed exactly at end-of-paragraph
# mis-places the bookmarkStart tag at the start of the paragraph instead).
# The placeholder is removed again once the bookmark has been created.
$r.InsertAfter(" 2X")
# Pull the range back in by one character so it sits right before the "X"
# placeholder (i.e. right after " 2"), then collapse it there.
[void]$r.MoveEnd(1, -1)
$r.Collapse(0)

# --- 2. Move the "_GoBack" bookmark to the new end of the text ----------
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
}

# $r is currently collapsed right before the "X" placeholder, i.e. right
# after " 2" -- exactly where the bookmark should end up.
$d.Bookmarks.Add("_GoBack", $r)

# --- 3. Remove the placeholder character ---------------------------------
$placeholder = $d.Range($r.End, $r.End + 1)
$placeholder.Delete()
